# Applies the FFXIV leve-profit refresh captured in the commit diff:
# recomputed currentAveragePrice* / LevePrice* / LeveProfit* columns (H, I, J, K, L, M, N)
# for a handful of rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 867.6429000000001
$ws.Range("J121").Value = 719
$ws.Range("L121").Value = 2157
$ws.Range("N121").Value = -5651
$ws.Range("H129").Value = 989.8148
$ws.Range("J129").Value = 1046.3673
$ws.Range("L129").Value = 3139.1019
$ws.Range("N129").Value = -13139.1019
$ws.Range("H134").Value = 67358
$ws.Range("J134").Value = 67358
$ws.Range("L134").Value = 67358
$ws.Range("N134").Value = -77498
$ws.Range("H140").Value = 109190
$ws.Range("J140").Value = 109190
$ws.Range("L140").Value = 109190
$ws.Range("N140").Value = -119550

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("H9").Value = 10009
$ws.Range("J9").Value = 10009
$ws.Range("L9").Value = 10009
$ws.Range("N9").Value = -10349
$ws.Range("H20").Value = 10009
$ws.Range("J20").Value = 10009
$ws.Range("L20").Value = 10009
$ws.Range("N20").Value = -10549
$ws.Range("H37").Value = 5034
$ws.Range("I37").Value = 5034
$ws.Range("K37").Value = 5034
$ws.Range("M37").Value = -4761
$ws.Range("H44").Value = 49666.668
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 69500
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 69500
$ws.Range("M44").Value = -9512
$ws.Range("N44").Value = -70476
$ws.Range("H55").Value = 39800
$ws.Range("I55").Value = 39800
$ws.Range("K55").Value = 39800
$ws.Range("M55").Value = -39485
$ws.Range("H80").Value = 40110
$ws.Range("J80").Value = 40110
$ws.Range("L80").Value = 40110
$ws.Range("N80").Value = -42106
$ws.Range("H83").Value = 40110
$ws.Range("J83").Value = 40110
$ws.Range("L83").Value = 120330
$ws.Range("N83").Value = -130314
$ws.Range("H129").Value = 36474.5
$ws.Range("J129").Value = 36474.5
$ws.Range("L129").Value = 36474.5
$ws.Range("N129").Value = -46474.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 70
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -296
$ws.Range("H31").Value = 715789.9
$ws.Range("I31").Value = 6618.3335
$ws.Range("J31").Value = 1091233.6
$ws.Range("K31").Value = 6618.3335
$ws.Range("L31").Value = 1091233.6
$ws.Range("M31").Value = -6323.3335
$ws.Range("N31").Value = -1091823.6
$ws.Range("H34").Value = 715789.9
$ws.Range("I34").Value = 6618.3335
$ws.Range("J34").Value = 1091233.6
$ws.Range("K34").Value = 6618.3335
$ws.Range("L34").Value = 1091233.6
$ws.Range("M34").Value = -6416.3335
$ws.Range("N34").Value = -1091637.6
$ws.Range("H92").Value = 40601
$ws.Range("J92").Value = 40601
$ws.Range("L92").Value = 40601
$ws.Range("N92").Value = -45593
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H140").Value = 82897.5
$ws.Range("J140").Value = 82897.5
$ws.Range("L140").Value = 82897.5
$ws.Range("N140").Value = -93257.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2446.197
$ws.Range("I68").Value = 827.3684
$ws.Range("J68").Value = 4643.1787
$ws.Range("K68").Value = 2482.1052
$ws.Range("L68").Value = 13929.5361
$ws.Range("M68").Value = -1671.1052
$ws.Range("N68").Value = -15551.5361
$ws.Range("H71").Value = 2446.197
$ws.Range("I71").Value = 827.3684
$ws.Range("J71").Value = 4643.1787
$ws.Range("K71").Value = 7446.3156
$ws.Range("L71").Value = 41788.60830000001
$ws.Range("M71").Value = -3390.3156
$ws.Range("N71").Value = -49900.60830000001
$ws.Range("H92").Value = 300
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = ""
$ws.Range("H94").Value = 8800
$ws.Range("J94").Value = 8800
$ws.Range("L94").Value = 26400
$ws.Range("N94").Value = -27752

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""
$ws.Range("H141").Value = 37457.25
$ws.Range("J141").Value = 37457.25
$ws.Range("L141").Value = 37457.25
$ws.Range("N141").Value = -47817.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2846
$ws.Range("I2").Value = 385
$ws.Range("J2").Value = 3666.3333
$ws.Range("K2").Value = 385
$ws.Range("L2").Value = 3666.3333
$ws.Range("M2").Value = -273
$ws.Range("N2").Value = -3890.3333
$ws.Range("H122").Value = 6036.8125
$ws.Range("I122").Value = 6549.1665
$ws.Range("J122").Value = 4499.75
$ws.Range("K122").Value = 19647.4995
$ws.Range("L122").Value = 13499.25
$ws.Range("M122").Value = -17197.4995
$ws.Range("N122").Value = -18399.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3573.5715
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 3573.5715
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H135").Value = 57460.332
$ws.Range("J135").Value = 57460.332
$ws.Range("L135").Value = 57460.332
$ws.Range("N135").Value = -67600.33199999999
$ws.Range("H140").Value = 41143
$ws.Range("J140").Value = 41143
$ws.Range("L140").Value = 41143
$ws.Range("N140").Value = -51503
$ws.Range("H141").Value = 45816.43
$ws.Range("J141").Value = 45816.43
$ws.Range("L141").Value = 45816.43
$ws.Range("N141").Value = -56176.43

Write-Output "Updated 159 cells across 8 sheets"
